$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Exact "17÷7=2, 3" "67÷9=7, 4"
Replace-Exact "58÷7=8, 2" "76÷6=12, 4"
Replace-Exact "65÷3=21, 2" "24÷4=6, 0"
Replace-Exact "24÷3=8, 0" "17÷4=4, 1"
Replace-Exact "55÷9=6, 1" "64÷4=16, 0"
Replace-Exact "24÷5=4, 4" "34÷7=4, 6"
Replace-Exact "52÷4=13, 0" "49÷2=24, 1"
Replace-Exact "64÷7=9, 1" "39÷6=6, 3"
Replace-Exact "45÷7=6, 3" "92÷7=13, 1"
Replace-Exact "75÷6=12, 3" "74÷5=14, 4"
Replace-Exact "79÷5=15, 4" "70÷7=10, 0"
Replace-Exact "76÷2=38, 0" "94÷7=13, 3"
Replace-Exact "54÷5=10, 4" "46÷6=7, 4"
Replace-Exact "44÷6=7, 2" "99÷2=49, 1"
Replace-Exact "35÷9=3, 8" "24÷7=3, 3"
Replace-Exact "93÷4=23, 1" "95÷2=47, 1"
Replace-Exact "10÷2=5, 0" "90÷6=15, 0"
Replace-Exact "86÷8=10, 6" "79÷2=39, 1"
Replace-Exact "43÷5=8, 3" "41÷6=6, 5"
Replace-Exact "60÷6=10, 0" "80÷7=11, 3"
Replace-Exact "17÷5=3, 2" "68÷5=13, 3"
Replace-Exact "42÷6=7, 0" "57÷8=7, 1"
Replace-Exact "56÷8=7, 0" "92÷7=13, 1"
Replace-Exact "54÷3=18, 0" "17÷8=2, 1"
Replace-Exact "78÷4=19, 2" "61÷2=30, 1"
